$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input Data")

# Remove the stored credentials (username/password values) while keeping
# the row labels in column A and the unrelated "OnePlus" search value.
$ws.Range("B1").ClearContents()
$ws.Range("B2").ClearContents()

$ws.Range("B1").Select()
